# Scheduled-runner market data refresh for the Ragnarok_Profits workbook.
# Re-prices each job sheet's currentAveragePrice* / LevePrice* / LeveProfit*
# columns (H:N) with freshly pulled Market Board data.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1189.5
$ws.Range("I2").Value = 773
$ws.Range("K2").Value = 773
$ws.Range("M2").Value = -660
$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("M34").ClearContents()
$ws.Range("H36").Value = 0
$ws.Range("I36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("M36").ClearContents()
$ws.Range("H53").Value = 1415
$ws.Range("J53").Value = 1749.25
$ws.Range("L53").Value = 1749.25
$ws.Range("N53").Value = -3023.25
$ws.Range("H62").Value = 5133.0454
$ws.Range("I62").Value = 4796.0586
$ws.Range("J62").Value = 6278.8
$ws.Range("K62").Value = 4796.0586
$ws.Range("L62").Value = 6278.8
$ws.Range("M62").Value = -4172.0586
$ws.Range("N62").Value = -7526.8
$ws.Range("H65").Value = 5133.0454
$ws.Range("I65").Value = 4796.0586
$ws.Range("J65").Value = 6278.8
$ws.Range("K65").Value = 23980.293
$ws.Range("L65").Value = 31394
$ws.Range("M65").Value = -20860.293
$ws.Range("N65").Value = -37634
$ws.Range("H76").Value = 13716.5
$ws.Range("I76").Value = 22988
$ws.Range("J76").Value = 4445
$ws.Range("K76").Value = 22988
$ws.Range("L76").Value = 4445
$ws.Range("M76").Value = -22673
$ws.Range("N76").Value = -5075
$ws.Range("H79").Value = 13716.5
$ws.Range("I79").Value = 22988
$ws.Range("J79").Value = 4445
$ws.Range("K79").Value = 22988
$ws.Range("L79").Value = 4445
$ws.Range("M79").Value = -21896
$ws.Range("N79").Value = -6629
$ws.Range("H100").Value = 6644
$ws.Range("J100").Value = 9889.200000000001
$ws.Range("L100").Value = 9889.200000000001
$ws.Range("N100").Value = -10971.2
$ws.Range("H103").Value = 29413038
$ws.Range("J103").Value = 50001644
$ws.Range("L103").Value = 150004932
$ws.Range("N103").Value = -150006104
$ws.Range("H137").Value = 2347.3845
$ws.Range("I137").Value = 1443.1333
$ws.Range("K137").Value = 4329.3999
$ws.Range("M137").Value = -1779.3999
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6533.356
$ws.Range("I32").Value = 5924.4463
$ws.Range("K32").Value = 5924.4463
$ws.Range("M32").Value = -5637.4463
$ws.Range("H110").Value = 5324.1055
$ws.Range("I110").Value = 4854.5
$ws.Range("J110").Value = 6639
$ws.Range("K110").Value = 4854.5
$ws.Range("L110").Value = 6639
$ws.Range("M110").Value = -2809.5
$ws.Range("N110").Value = -10729
$ws.Range("H122").Value = 4374.5
$ws.Range("I122").Value = 4879.5
$ws.Range("J122").Value = 1849.5
$ws.Range("K122").Value = 14638.5
$ws.Range("L122").Value = 5548.5
$ws.Range("M122").Value = -12188.5
$ws.Range("N122").Value = -10448.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2406.5
$ws.Range("I94").Value = 3329.7273
$ws.Range("J94").Value = 955.7143
$ws.Range("K94").Value = 3329.7273
$ws.Range("L94").Value = 955.7143
$ws.Range("M94").Value = -2878.7273
$ws.Range("N94").Value = -1857.7143
$ws.Range("H134").Value = 3574198
$ws.Range("I134").Value = 2759.652
$ws.Range("K134").Value = 8278.956
$ws.Range("M134").Value = -5743.956

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 45457936
$ws.Range("I31").Value = 58826224
$ws.Range("K31").Value = 58826224
$ws.Range("M31").Value = -58825929
$ws.Range("H34").Value = 45457936
$ws.Range("I34").Value = 58826224
$ws.Range("K34").Value = 58826224
$ws.Range("M34").Value = -58826022
$ws.Range("H58").Value = 4433
$ws.Range("I58").Value = 3648.4
$ws.Range("K58").Value = 3648.4
$ws.Range("M58").Value = -3445.4
$ws.Range("H86").Value = 7993.353
$ws.Range("I86").Value = 6207.154
$ws.Range("K86").Value = 6207.154
$ws.Range("M86").Value = -5084.154
$ws.Range("H89").Value = 7993.353
$ws.Range("I89").Value = 6207.154
$ws.Range("K89").Value = 31035.77
$ws.Range("M89").Value = -25419.77
$ws.Range("H132").Value = 1841
$ws.Range("I132").Value = 1515.4117
$ws.Range("J132").Value = 3224.75
$ws.Range("K132").Value = 4546.2351
$ws.Range("L132").Value = 9674.25
$ws.Range("M132").Value = -2016.2351
$ws.Range("N132").Value = -14734.25
$ws.Range("H136").Value = 4433
$ws.Range("I136").Value = 3648.4
$ws.Range("K136").Value = 10945.2
$ws.Range("M136").Value = -8395.200000000001
$ws.Range("H137").Value = 55211
$ws.Range("I137").Value = 55211
$ws.Range("K137").Value = 55211
$ws.Range("M137").Value = -50111

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H117").Value = 5187.6
$ws.Range("I117").Value = 483.66666
$ws.Range("J117").Value = 7203.5713
$ws.Range("K117").Value = 1450.99998
$ws.Range("L117").Value = 21610.7139
$ws.Range("M117").Value = 1991.00002
$ws.Range("N117").Value = -28494.7139
$ws.Range("H122").Value = 66282.8
$ws.Range("I122").Value = 82728.5
$ws.Range("K122").Value = 744556.5
$ws.Range("M122").Value = -742106.5
$ws.Range("H133").Value = 49673.668
$ws.Range("I133").Value = 52399.8
$ws.Range("J133").Value = 36043
$ws.Range("K133").Value = 157199.4
$ws.Range("L133").Value = 108129
$ws.Range("M133").Value = -152139.4
$ws.Range("N133").Value = -118249

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 11573.721
$ws.Range("J70").Value = 11323.182
$ws.Range("L70").Value = 11323.182
$ws.Range("N70").Value = -11863.182
$ws.Range("H73").Value = 11573.721
$ws.Range("J73").Value = 11323.182
$ws.Range("L73").Value = 11323.182
$ws.Range("N73").Value = -13195.182
$ws.Range("H132").Value = 2781435
$ws.Range("I132").Value = 3744.0645
$ws.Range("K132").Value = 11232.1935
$ws.Range("M132").Value = -8702.193499999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8268.462
$ws.Range("I7").Value = 8268.462
$ws.Range("K7").Value = 8268.462
$ws.Range("M7").Value = -8156.462
$ws.Range("H93").Value = 2927434.2
$ws.Range("I93").Value = 3098.4443
$ws.Range("K93").Value = 3098.4443
$ws.Range("M93").Value = -1850.4443
$ws.Range("H116").Value = 117499.5
$ws.Range("J116").Value = 117499.5
$ws.Range("L116").Value = 117499.5
$ws.Range("N116").Value = -126677.5
$ws.Range("H122").Value = 4786.9766
$ws.Range("I122").Value = 3495.5
$ws.Range("K122").Value = 10486.5
$ws.Range("M122").Value = -8036.5
$ws.Range("H126").Value = 8268.462
$ws.Range("I126").Value = 8268.462
$ws.Range("K126").Value = 24805.386
$ws.Range("M126").Value = -22335.386

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 37499.5
$ws.Range("J2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("N2").ClearContents()
$ws.Range("H44").Value = 33941
$ws.Range("J44").Value = 33941
$ws.Range("L44").Value = 33941
$ws.Range("N44").Value = -35049
$ws.Range("H45").Value = 22445.75
$ws.Range("I45").Value = 30899.5
$ws.Range("J45").Value = 13992
$ws.Range("K45").Value = 30899.5
$ws.Range("L45").Value = 13992
$ws.Range("M45").Value = -30408.5
$ws.Range("N45").Value = -14974
$ws.Range("H46").Value = 78962
$ws.Range("J46").Value = 78962
$ws.Range("L46").Value = 78962
$ws.Range("N46").Value = -79424
$ws.Range("H49").Value = 33962
$ws.Range("J49").Value = 33962
$ws.Range("L49").Value = 33962
$ws.Range("N49").Value = -34422
$ws.Range("H81").Value = 2084
$ws.Range("I81").Value = 2159.625
$ws.Range("K81").Value = 4319.25
$ws.Range("M81").Value = -3258.25
$ws.Range("H84").Value = 2084
$ws.Range("I84").Value = 2159.625
$ws.Range("K84").Value = 21596.25
$ws.Range("M84").Value = -16292.25
$ws.Range("H122").Value = 2084.5625
$ws.Range("I122").Value = 1668.6364
$ws.Range("J122").Value = 2999.6
$ws.Range("K122").Value = 5005.9092
$ws.Range("L122").Value = 8998.799999999999
$ws.Range("M122").Value = -2555.9092
$ws.Range("N122").Value = -13898.8
$ws.Range("H126").Value = 4343.222
$ws.Range("I126").Value = 5608
$ws.Range("J126").Value = 2762.25
$ws.Range("K126").Value = 16824
$ws.Range("L126").Value = 8286.75
$ws.Range("M126").Value = -14354
$ws.Range("N126").Value = -13226.75
$ws.Range("H132").Value = 224603.33
$ws.Range("I132").Value = 2496.5264
$ws.Range("J132").Value = 1430326
$ws.Range("K132").Value = 7489.5792
$ws.Range("L132").Value = 4290978
$ws.Range("M132").Value = -4959.5792
$ws.Range("N132").Value = -4296038
$ws.Range("H134").Value = 78962
$ws.Range("J134").Value = 78962
$ws.Range("L134").Value = 236886
$ws.Range("N134").Value = -241956
$ws.Range("H141").Value = 89843
$ws.Range("J141").Value = 89843
$ws.Range("L141").Value = 89843
$ws.Range("N141").Value = -100203
